$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Parametric"
$ws.Range("C13").Select()
